# Generate Report for Handback
#
# The last handback batch (row 16, file "de8179a3-...") in each language
# sheet had its "Correspond Handoff Datetime" (col D) / "Correspond
# Handback DateTime" (col G) cells still showing the placeholder values
# copied down from the row above (row 15) while the real handback was in
# flight. Now that the handback for that file has completed, stamp the
# real timestamps recorded for it onto row 16 of each language sheet.

$wb = $excel.ActiveWorkbook

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_zhcn.Range("D16").Value = "2016-03-09 04:58:49"
$ws_zhcn.Range("G16").Value = "2016-03-09 04:59:45"

$ws_dede = $wb.Worksheets.Item("de-de")
$ws_dede.Range("D16").Value = "2016-03-09 04:58:51"
$ws_dede.Range("G16").Value = "2016-03-09 04:59:50"
